$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "67.204.53"
$ws.Range("E2").Value = "  +5.19%  "
Set-TextValue "D3" "3.517.58"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "594.95"
$ws.Range("E5").Value = "  +4.24%  "
Set-TextValue "D6" "169.37"
$ws.Range("E6").Value = "  +7.14%  "
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.05%  "
Set-TextValue "D8" "3.519.72"
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("E11").Value = "  +5.92%  "
$ws.Range("E12").Value = "  +4.79%  "
Set-TextValue "D13" "4.123.92"
$ws.Range("E13").Value = "  +2.96%  "
Set-TextValue "D15" "28.31"
$ws.Range("E15").Value = "  +4.60%  "
$ws.Range("E16").Value = "  +4.59%  "
Set-TextValue "D17" "67.184.73"
$ws.Range("E17").Value = "  +5.11%  "
Set-TextValue "D18" "3.508.43"
$ws.Range("E18").Value = "  +2.26%  "
Set-TextValue "D19" "6.34"
$ws.Range("E19").Value = "  +4.27%  "
Set-TextValue "D20" "14.08"
$ws.Range("E20").Value = "  +3.47%  "
Set-TextValue "D21" "396.28"
$ws.Range("E21").Value = "  +3.41%  "
Set-TextValue "D22" "7.97"
$ws.Range("E22").Value = "  +2.27%  "
Set-TextValue "D23" "73.61"
$ws.Range("E23").Value = "  +3.26%  "
$ws.Range("E24").Value = "  +11.25%  "
Set-TextValue "D25" "0.999"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  +3.44%  "
Set-TextValue "D27" "10.22"
$ws.Range("E27").Value = "  +5.58%  "
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +5.49%  "
$ws.Range("E31").Value = "  +6.48%  "
$ws.Range("E32").Value = "  +4.18%  "
Set-TextValue "D33" "23.66"
$ws.Range("E33").Value = "  +3.38%  "
Set-TextValue "D34" "7.47"
$ws.Range("E34").Value = "  +7.49%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  +6.37%  "
Set-TextValue "D37" "161.90"
$ws.Range("E37").Value = "  +0.68%  "
Set-TextValue "D38" "0.903"
$ws.Range("E38").Value = "  +6.61%  "
Set-TextValue "D39" "1.94"
$ws.Range("E39").Value = "  +6.61%  "
Set-TextValue "D40" "0.0754"
$ws.Range("E40").Value = "  +4.53%  "
$ws.Range("E41").Value = "  +7.39%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D42" "26.66"
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D43" "6.72"
$ws.Range("E43").Value = "  +5.30%  "
Set-TextValue "D44" "2.845.32"
$ws.Range("E44").Value = "  +1.76%  "
Set-TextValue "D45" "43.52"
$ws.Range("E45").Value = "  +1.06%  "
Set-TextValue "D46" "26.60"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D47" "0.0316"
$ws.Range("E47").Value = "  +4.01%  "
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D48" "2.56"
$ws.Range("E48").Value = "  +6.72%  "
Set-TextValue "D49" "355.19"
$ws.Range("E49").Value = "  +6.78%  "
$ws.Range("E50").Value = "  +5.09%  "
Set-TextValue "D51" "33.70"
$ws.Range("E51").Value = "  +12.88%  "
